# Revision: add "chemical_recycling_pyrolysis" as a new boolean parameter
# right after "chemical_recycling_gasification" (row 9), shifting every
# subsequent parameter row down by one. All other rows keep their original
# values/explanations - only their row position changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 10; rows 10-24 shift down to 11-25, and the
# used range grows from A1:C24 to A1:C25.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row with the new parameter.
$ws.Range("A10").Value = "chemical_recycling_pyrolysis"
$ws.Range("B10").Value = $true
